# Applies the cryptos-list refresh described in the commit diff:
# updated Price/Volume(1h) figures for every row, plus a rank swap
# between Filecoin/HuobiToken (rows 34/36) and Hedera/TrustWalletToken
# (rows 38/39).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '28.658.73'; E = '  -3.20%  ' },
    @{ Row = 3; D = '1.955.29'; E = '  -2.42%  ' },
    @{ Row = 4; D = '1.015'; E = '  +0.16%  ' },
    @{ Row = 5; D = '321.64'; E = '  -2.41%  ' },
    @{ Row = 6; D = '1.013'; E = '  +0.20%  ' },
    @{ Row = 7; D = '0.4764'; E = '  -4.85%  ' },
    @{ Row = 8; D = '0.4035'; E = '  -4.66%  ' },
    @{ Row = 9; D = '53.74'; E = '  -0.46%  ' },
    @{ Row = 10; D = '0.08402'; E = '  -6.85%  ' },
    @{ Row = 11; D = '1.060'; E = '  -5.30%  ' },
    @{ Row = 12; D = '22.18'; E = '  -5.07%  ' },
    @{ Row = 13; D = '1.935.26'; E = '  -2.78%  ' },
    @{ Row = 14; D = '7.599'; E = '  -5.79%  ' },
    @{ Row = 15; D = '6.193'; E = '  -4.42%  ' },
    @{ Row = 16; D = '1.015'; E = '  +0.22%  ' },
    @{ Row = 17; D = '0.00001076'; E = '  -3.48%  ' },
    @{ Row = 18; D = '89.02'; E = '  -5.51%  ' },
    @{ Row = 19; D = '0.06639'; E = '  -0.41%  ' },
    @{ Row = 20; D = '18.66'; E = '  -5.34%  ' },
    @{ Row = 21; D = '1.013'; E = '  +0.11%  ' },
    @{ Row = 22; D = '5.815'; E = '  -2.35%  ' },
    @{ Row = 23; D = '28.680.66'; E = '  -3.15%  ' },
    @{ Row = 24; D = '11.52'; E = '  -3.84%  ' },
    @{ Row = 25; D = '2.298'; E = '  -0.15%  ' },
    @{ Row = 26; D = '2.172.58'; E = '  -2.28%  ' },
    @{ Row = 27; D = '154.33'; E = '  -2.64%  ' },
    @{ Row = 28; D = '20.14'; E = '  -2.78%  ' },
    @{ Row = 29; D = '5.918'; E = '  -7.73%  ' },
    @{ Row = 30; D = '2.154'; E = '  -6.56%  ' },
    @{ Row = 31; D = '123.46'; E = '  -3.79%  ' },
    @{ Row = 32; D = '0.9993'; E = '  -5.19%  ' },
    @{ Row = 33; D = '0.09587'; E = '  -3.58%  ' },
    @{ Row = 34; B = 'HuobiToken'; C = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D = '3.668'; E = '  -3.44%  ' },
    @{ Row = 35; D = '1.430'; E = '  -9.22%  ' },
    @{ Row = 36; B = 'Filecoin'; C = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D = '5.565'; E = '  -4.73%  ' },
    @{ Row = 37; D = '0.02331'; E = '  -5.76%  ' },
    @{ Row = 38; B = 'Hedera'; C = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; D = '0.06226'; E = '  -2.11%  ' },
    @{ Row = 39; B = 'TrustWalletToken'; C = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D = '1.260'; E = '  -3.68%  ' },
    @{ Row = 40; D = '8.724'; E = '  -6.33%  ' },
    @{ Row = 41; D = '0.6217'; E = '  -5.29%  ' },
    @{ Row = 42; D = '11.07'; E = '  -5.31%  ' },
    @{ Row = 43; D = '1.012'; E = '  +0.18%  ' },
    @{ Row = 44; D = '0.1918'; E = '  -6.51%  ' },
    @{ Row = 45; D = '1.329'; E = '  +1.88%  ' },
    @{ Row = 46; D = '0.5936'; E = '  -6.57%  ' },
    @{ Row = 47; D = '12.92'; E = '  -3.73%  ' },
    @{ Row = 48; D = '2.070'; E = '  -5.78%  ' },
    @{ Row = 49; D = '3.416'; E = '  -2.45%  ' },
    @{ Row = 50; D = '0.00000000334'; E = '  +0.01%  ' },
    @{ Row = 51; D = '0.06824'; E = '  -2.34%  ' }
)

foreach ($item in $updates) {
    $r = $item.Row
    if ($item.ContainsKey('B')) { $ws.Range("B$r").Value = $item.B }
    if ($item.ContainsKey('C')) { $ws.Range("C$r").Value = $item.C }
    if ($item.ContainsKey('D')) {
        $dCell = $ws.Range("D$r")
        $dVal = $item.D
        # The Price column stores plain text (e.g. "1.015", "28.658.73").
        # Force text formatting first so values that parse as plain numbers
        # aren't silently coerced to doubles by the COM Value setter --
        # this mirrors how the sheet's Price cells are authored upstream.
        if ($dVal -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $dVal
    }
    if ($item.ContainsKey('E')) { $ws.Range("E$r").Value = $item.E }
}